$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy style s="1" (bold + border + centered) from an existing labeled cell
# into the two brand-new J-column label cells (J20, J21) before assigning values.
$ws.Range("J19").Copy()
$ws.Range("J20:J21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A1").Value = 'negative'
$ws.Range("J1").Value = 'positive'
$ws.Range("A2").Value = 'name'
$ws.Range("B2").Value = 'anchor score'
$ws.Range("C2").Value = 'type occurences'
$ws.Range("D2").Value = 'total occurences'
$ws.Range("E2").Value = '+%'
$ws.Range("F2").Value = '-%'
$ws.Range("G2").Value = 'both'
$ws.Range("H2").Value = 'normal'
$ws.Range("J2").Value = 'name'
$ws.Range("K2").Value = 'anchor score'
$ws.Range("L2").Value = 'type occurences'
$ws.Range("M2").Value = 'total occurences'
$ws.Range("N2").Value = '+%'
$ws.Range("O2").Value = '-%'
$ws.Range("P2").Value = 'both'
$ws.Range("Q2").Value = 'normal'
$ws.Range("A3").Value = 'poorly'
$ws.Range("B3").Value = 0.9782608695652174
$ws.Range("C3").Value = 45
$ws.Range("D3").Value = 45
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = $false
$ws.Range("H3").Value = 1
$ws.Range("J3").Value = 'wonderful'
$ws.Range("K3").Value = 0.8928571428571429
$ws.Range("L3").Value = 50
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = $false
$ws.Range("Q3").Value = 6
$ws.Range("A4").Value = 'disappointing'
$ws.Range("B4").Value = 0.9090909090909091
$ws.Range("C4").Value = 40
$ws.Range("D4").Value = 40
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = 4
$ws.Range("J4").Value = 'awesome'
$ws.Range("K4").Value = 0.8923076923076924
$ws.Range("L4").Value = 58
$ws.Range("M4").Value = 58
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = $false
$ws.Range("Q4").Value = 7
$ws.Range("A5").Value = 'poor'
$ws.Range("B5").Value = 0.7464788732394366
$ws.Range("C5").Value = 53
$ws.Range("D5").Value = 53
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = 18
$ws.Range("J5").Value = 'favorite'
$ws.Range("K5").Value = 0.8172043010752689
$ws.Range("L5").Value = 76
$ws.Range("M5").Value = 76
$ws.Range("N5").Value = 1
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = $false
$ws.Range("Q5").Value = 17
$ws.Range("A6").Value = 'however'
$ws.Range("B6").Value = 0.71875
$ws.Range("C6").Value = 46
$ws.Range("D6").Value = 46
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = $false
$ws.Range("H6").Value = 18
$ws.Range("J6").Value = 'excellent'
$ws.Range("K6").Value = 0.796875
$ws.Range("L6").Value = 51
$ws.Range("M6").Value = 51
$ws.Range("N6").Value = 1
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = $false
$ws.Range("Q6").Value = 13
$ws.Range("A7").Value = 'disappointed'
$ws.Range("B7").Value = 0.6720430107526881
$ws.Range("C7").Value = 125
$ws.Range("D7").Value = 125
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = $false
$ws.Range("H7").Value = 61
$ws.Range("J7").Value = 'classic'
$ws.Range("K7").Value = 0.6792452830188679
$ws.Range("L7").Value = 36
$ws.Range("M7").Value = 36
$ws.Range("N7").Value = 1
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = $false
$ws.Range("Q7").Value = 17
$ws.Range("A8").Value = 'broke'
$ws.Range("B8").Value = 0.6359223300970874
$ws.Range("C8").Value = 131
$ws.Range("D8").Value = 131
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = $false
$ws.Range("H8").Value = 75
$ws.Range("J8").Value = 'love'
$ws.Range("K8").Value = 0.5538020086083214
$ws.Range("L8").Value = 386
$ws.Range("M8").Value = 386
$ws.Range("N8").Value = 1
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = $false
$ws.Range("Q8").Value = 311
$ws.Range("A9").Value = 'waste'
$ws.Range("B9").Value = 0.6283783783783784
$ws.Range("C9").Value = 93
$ws.Range("D9").Value = 93
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = $false
$ws.Range("H9").Value = 55
$ws.Range("J9").Value = 'thank'
$ws.Range("K9").Value = 0.5362318840579711
$ws.Range("L9").Value = 37
$ws.Range("M9").Value = 37
$ws.Range("N9").Value = 1
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = $false
$ws.Range("Q9").Value = 32
$ws.Range("A10").Value = 'smaller'
$ws.Range("B10").Value = 0.5798319327731093
$ws.Range("C10").Value = 69
$ws.Range("D10").Value = 69
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = $false
$ws.Range("H10").Value = 50
$ws.Range("J10").Value = 'loves'
$ws.Range("K10").Value = 0.495850622406639
$ws.Range("L10").Value = 239
$ws.Range("M10").Value = 239
$ws.Range("N10").Value = 1
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = $false
$ws.Range("Q10").Value = 243
$ws.Range("A11").Value = 'junk'
$ws.Range("B11").Value = 0.5454545454545454
$ws.Range("C11").Value = 30
$ws.Range("D11").Value = 30
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = $false
$ws.Range("H11").Value = 25
$ws.Range("J11").Value = 'great'
$ws.Range("K11").Value = 0.4778688524590164
$ws.Range("L11").Value = 583
$ws.Range("M11").Value = 583
$ws.Range("N11").Value = 1
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = $false
$ws.Range("Q11").Value = 637
$ws.Range("A12").Value = 'small'
$ws.Range("B12").Value = 0.4898550724637681
$ws.Range("C12").Value = 169
$ws.Range("D12").Value = 169
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = $false
$ws.Range("H12").Value = 176
$ws.Range("J12").Value = 'loved'
$ws.Range("K12").Value = 0.345565749235474
$ws.Range("L12").Value = 113
$ws.Range("M12").Value = 113
$ws.Range("N12").Value = 1
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = $false
$ws.Range("Q12").Value = 214
$ws.Range("A13").Value = 'broken'
$ws.Range("B13").Value = 0.4337349397590362
$ws.Range("C13").Value = 36
$ws.Range("D13").Value = 36
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = $false
$ws.Range("H13").Value = 47
$ws.Range("J13").Value = 'friends'
$ws.Range("K13").Value = 0.3227513227513227
$ws.Range("L13").Value = 61
$ws.Range("M13").Value = 61
$ws.Range("N13").Value = 1
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = $false
$ws.Range("Q13").Value = 128
$ws.Range("A14").Value = 'plastic'
$ws.Range("B14").Value = 0.3937007874015748
$ws.Range("C14").Value = 50
$ws.Range("D14").Value = 50
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = $false
$ws.Range("H14").Value = 77
$ws.Range("J14").Value = 'perfect'
$ws.Range("K14").Value = 0.3192771084337349
$ws.Range("L14").Value = 53
$ws.Range("M14").Value = 53
$ws.Range("N14").Value = 1
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = $false
$ws.Range("Q14").Value = 113
$ws.Range("A15").Value = 'cheap'
$ws.Range("B15").Value = 0.3649289099526066
$ws.Range("C15").Value = 77
$ws.Range("D15").Value = 77
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = $false
$ws.Range("H15").Value = 134
$ws.Range("J15").Value = 'best'
$ws.Range("K15").Value = 0.3166666666666667
$ws.Range("L15").Value = 38
$ws.Range("M15").Value = 38
$ws.Range("N15").Value = 1
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = $false
$ws.Range("Q15").Value = 82
$ws.Range("A16").Value = 'apart'
$ws.Range("B16").Value = 0.3473684210526316
$ws.Range("C16").Value = 33
$ws.Range("D16").Value = 33
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = $false
$ws.Range("H16").Value = 62
$ws.Range("J16").Value = 'happy'
$ws.Range("K16").Value = 0.2097902097902098
$ws.Range("L16").Value = 30
$ws.Range("M16").Value = 30
$ws.Range("N16").Value = 1
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = $false
$ws.Range("Q16").Value = 113
$ws.Range("A17").Value = 'difficult'
$ws.Range("B17").Value = 0.3370786516853932
$ws.Range("C17").Value = 30
$ws.Range("D17").Value = 30
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = $false
$ws.Range("H17").Value = 59
$ws.Range("J17").Value = 'enjoy'
$ws.Range("K17").Value = 0.1989247311827957
$ws.Range("L17").Value = 37
$ws.Range("M17").Value = 37
$ws.Range("N17").Value = 1
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = $false
$ws.Range("Q17").Value = 149
$ws.Range("A18").Value = 'ok'
$ws.Range("B18").Value = 0.3203125
$ws.Range("C18").Value = 41
$ws.Range("D18").Value = 41
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = $false
$ws.Range("H18").Value = 87
$ws.Range("J18").Value = 'christmas'
$ws.Range("K18").Value = 0.1566265060240964
$ws.Range("L18").Value = 39
$ws.Range("M18").Value = 39
$ws.Range("N18").Value = 1
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = $false
$ws.Range("Q18").Value = 210
$ws.Range("A19").Value = 'thought'
$ws.Range("B19").Value = 0.2574257425742574
$ws.Range("C19").Value = 52
$ws.Range("D19").Value = 52
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = $false
$ws.Range("H19").Value = 150
$ws.Range("J19").Value = 'fun'
$ws.Range("K19").Value = 0.1489921121822962
$ws.Range("L19").Value = 170
$ws.Range("M19").Value = 170
$ws.Range("N19").Value = 1
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = $false
$ws.Range("Q19").Value = 971
$ws.Range("A20").Value = 'size'
$ws.Range("B20").Value = 0.2268041237113402
$ws.Range("C20").Value = 44
$ws.Range("D20").Value = 44
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = $false
$ws.Range("H20").Value = 150
$ws.Range("J20").Value = 'family'
$ws.Range("K20").Value = 0.08913649025069638
$ws.Range("L20").Value = 32
$ws.Range("M20").Value = 32
$ws.Range("N20").Value = 1
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = $false
$ws.Range("Q20").Value = 327
$ws.Range("A21").Value = 'item'
$ws.Range("B21").Value = 0.1884057971014493
$ws.Range("C21").Value = 52
$ws.Range("D21").Value = 52
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = $false
$ws.Range("H21").Value = 224
$ws.Range("J21").Value = 'game'
$ws.Range("K21").Value = 0.07268007787151201
$ws.Range("L21").Value = 112
$ws.Range("M21").Value = 112
$ws.Range("N21").Value = 1
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = $false
$ws.Range("Q21").Value = 1429
$ws.Range("A22").Value = 'money'
$ws.Range("B22").Value = 0.1740506329113924
$ws.Range("C22").Value = 55
$ws.Range("D22").Value = 55
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = $false
$ws.Range("H22").Value = 261
$ws.Range("A23").Value = 'better'
$ws.Range("B23").Value = 0.1728971962616822
$ws.Range("C23").Value = 37
$ws.Range("D23").Value = 37
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = $false
$ws.Range("H23").Value = 177
$ws.Range("A24").Value = 'hard'
$ws.Range("B24").Value = 0.17
$ws.Range("C24").Value = 34
$ws.Range("D24").Value = 34
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = $false
$ws.Range("H24").Value = 166
$ws.Range("A25").Value = 'would'
$ws.Range("B25").Value = 0.1691394658753709
$ws.Range("C25").Value = 114
$ws.Range("D25").Value = 114
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = $false
$ws.Range("H25").Value = 560
$ws.Range("A26").Value = 'price'
$ws.Range("B26").Value = 0.1609195402298851
$ws.Range("C26").Value = 56
$ws.Range("D26").Value = 56
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = $false
$ws.Range("H26").Value = 292
$ws.Range("A27").Value = 'work'
$ws.Range("B27").Value = 0.1582278481012658
$ws.Range("C27").Value = 50
$ws.Range("D27").Value = 50
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = $false
$ws.Range("H27").Value = 266
$ws.Range("A28").Value = 'product'
$ws.Range("B28").Value = 0.1167400881057269
$ws.Range("C28").Value = 53
$ws.Range("D28").Value = 53
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = $false
$ws.Range("H28").Value = 401
$ws.Range("A29").Value = '2'
$ws.Range("B29").Value = 0.1123595505617977
$ws.Range("C29").Value = 30
$ws.Range("D29").Value = 30
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = $false
$ws.Range("H29").Value = 237
$ws.Range("A30").Value = 'little'
$ws.Range("B30").Value = 0.06904231625835189
$ws.Range("C30").Value = 31
$ws.Range("D30").Value = 31
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = $false
$ws.Range("H30").Value = 418
$ws.Range("A31").Value = 'like'
$ws.Range("B31").Value = 0.05107084019769358
$ws.Range("C31").Value = 31
$ws.Range("D31").Value = 32
$ws.Range("E31").Value = 0.03
$ws.Range("F31").Value = 0.97
$ws.Range("G31").Value = $true
$ws.Range("H31").Value = 576
$ws.Range("A32").Value = 'one'
$ws.Range("B32").Value = 0.04562737642585551
$ws.Range("C32").Value = 36
$ws.Range("D32").Value = 41
$ws.Range("E32").Value = 0.12
$ws.Range("F32").Value = 0.88
$ws.Range("G32").Value = $true
$ws.Range("H32").Value = 753
